$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7..50 shift down to 8..51
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new market-day record
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Femacal de La Calera"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44749
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 100112035
$ws.Range("G7").Value = "Bruselas (repollito)"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 16450
$ws.Range("N7").Value = "`$/malla 15 kilos"
$ws.Range("O7").Value = "Provincia de Quillota"
$ws.Range("P7").Value = 1097
$ws.Range("Q7").Value = 15
$ws.Range("R7").Value = "Hortaliza"
